# update ANR MC scenario and MIBHolder
#
# 1) sheet ENODEB_R15_00 / ENODEB_R15_10: add the two ANR "nghRemoveThreshold"
#    rows (Cell1 / Cell2) that already exist on the other ENODEB_* tabs.
# 2) sheet ENODEB_R15_50 / ENODEB_R16_00 / ENODEB_R16_50 / ENODEB_R17_00:
#    the OID for those two ANR rows loses its leading "." prefix.
# 3) Misc cosmetic touch-ups: selection/tab-selection, ENODEB_R16_50 column A
#    width + row 52 height.

$wb = $excel.ActiveWorkbook

$shRev00 = $wb.Worksheets.Item("ENODEB_R15_00")
$shRev10 = $wb.Worksheets.Item("ENODEB_R15_10")
$shRev50 = $wb.Worksheets.Item("ENODEB_R15_50")
$shRev1600 = $wb.Worksheets.Item("ENODEB_R16_00")
$shRev1650 = $wb.Worksheets.Item("ENODEB_R16_50")
$shRev1700 = $wb.Worksheets.Item("ENODEB_R17_00")

# ---------------------------------------------------------------------------
# 1) Add the two missing ANR rows to ENODEB_R15_00 and ENODEB_R15_10.
#    Re-use the formatting already present on ENODEB_R15_50's equivalent
#    rows so the new cells look exactly like their siblings elsewhere in
#    the workbook.
# ---------------------------------------------------------------------------
foreach ($sh in @($shRev00, $shRev10)) {
    $sh.Range("A52").Value = "asLteStkCellAnrCfgNghRemoveThresholdCell1"
    $sh.Range("A52").Style = $shRev50.Range("A52").Style
    $sh.Range("B52").Value = "1.3.6.1.4.1.989.1.20.1.4.23.1.33.40"
    $sh.Range("B52").Style = $shRev50.Range("B52").Style
    $sh.Range("C52").Value = "nghRemoveThreshold"
    $sh.Range("C52").Style = $shRev50.Range("C52").Style
    $sh.Range("D52").Value = "Int"
    $sh.Range("D52").Style = $shRev50.Range("D52").Style

    $sh.Range("A53").Value = "asLteStkCellAnrCfgNghRemoveThresholdCell2"
    $sh.Range("A53").Style = $shRev50.Range("A53").Style
    $sh.Range("B53").Value = "1.3.6.1.4.1.989.1.20.1.4.23.1.33.41"
    $sh.Range("B53").Style = $shRev50.Range("B53").Style
    $sh.Range("C53").Value = "nghRemoveThreshold"
    $sh.Range("C53").Style = $shRev50.Range("C53").Style
    $sh.Range("D53").Value = "Int"
    $sh.Range("D53").Style = $shRev50.Range("D53").Style
}

# ---------------------------------------------------------------------------
# 2) Strip the leading "." from the ANR OIDs that already existed.
# ---------------------------------------------------------------------------
$shRev50.Range("B52").Value = "1.3.6.1.4.1.989.1.20.1.4.23.1.33.40"
$shRev50.Range("B53").Value = "1.3.6.1.4.1.989.1.20.1.4.23.1.33.41"

$shRev1600.Range("B52").Value = "1.3.6.1.4.1.989.1.20.1.4.23.1.33.40"
$shRev1600.Range("B53").Value = "1.3.6.1.4.1.989.1.20.1.4.23.1.33.41"

$shRev1650.Range("B53").Value = "1.3.6.1.4.1.989.1.20.1.4.23.1.33.40"
$shRev1650.Range("B54").Value = "1.3.6.1.4.1.989.1.20.1.4.23.1.33.41"

$shRev1700.Range("B52").Value = "1.3.6.1.4.1.989.1.20.1.4.23.1.33.40"
$shRev1700.Range("B53").Value = "1.3.6.1.4.1.989.1.20.1.4.23.1.33.41"

# ---------------------------------------------------------------------------
# 3) Cosmetic tweaks.
# ---------------------------------------------------------------------------

# ENODEB_R16_50: widen column A (now that it no longer auto-fits) and give
# row 52 its slightly reduced custom height.
$shRev1650.Columns.Item(1).ColumnWidth = 49.140625
$shRev1650.Rows.Item(52).RowHeight = 11.25

# Selection / active-tab bookkeeping: ENODEB_R17_00 is no longer the tab
# that's left selected - ENODEB_R15_00 is. Re-select the two new rows on
# each touched sheet so the saved selection rectangle matches, finishing on
# ENODEB_R15_00 so it becomes the active tab.
$shRev1700.Range("A52:XFD53").Select()
$shRev10.Range("A52:XFD53").Select()
$shRev00.Range("A52:XFD53").Select()
